# Update countries & provincias Spain
#
# Sheet "Pais" is a COVID-19 stats table (columns A:H = Pais, Casos
# totales, Nuevos casos, Casos activos, Recuperados, Casos criticos,
# Muertes hoy, Muertes) kept sorted by "Casos totales" (col B)
# descending. This refresh pulls newer per-country figures; a handful of
# countries leapfrog their neighbours in the sort as a result, which is
# why some rows below get a new country name (col A) in addition to new
# numbers, while the row right next to them effectively inherits the
# old numbers of the country that used to sit there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Range("A1").Value = 'Datos actualizados a 21 de Junio de 2020 a las 15:13'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2331550
$ws.Range("C4").Value = 972
$ws.Range("E4").Value = 1236492
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 122003

# India (row 7)
$ws.Range("B7").Value = 413092
$ws.Range("C7").Value = 1365
$ws.Range("D7").Value = 229148
$ws.Range("E7").Value = 170650
$ws.Range("G7").Value = 17
$ws.Range("H7").Value = 13294

# Arabia Saudita (row 19)
$ws.Range("B19").Value = 157612
$ws.Range("C19").Value = 3379
$ws.Range("D19").Value = 101130
$ws.Range("E19").Value = 55215
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 1267

# Paises Bajos (row 31)
$ws.Range("B31").Value = 49593
$ws.Range("C31").Value = 91
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 6090

# Portugal (row 37)
$ws.Range("B37").Value = 39133
$ws.Range("C37").Value = 292
$ws.Range("D37").Value = 25376
$ws.Range("E37").Value = 12227
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = 1530

# Uzbekistan (row 76)
$ws.Range("B76").Value = 6272
$ws.Range("C76").Value = 119
$ws.Range("E76").Value = 1963

# Tayikistan (row 79)
$ws.Range("B79").Value = 5457
$ws.Range("C79").Value = 58
$ws.Range("D79").Value = 3995
$ws.Range("E79").Value = 1410

# Kenia / Etiopia swap places (rows 85-86): Etiopia overtakes Kenia
$ws.Range("A85").Value = 'Etiopia'
$ws.Range("B85").Value = 4532
$ws.Range("C85").Value = 63
$ws.Range("D85").Value = 1213
$ws.Range("E85").Value = 3245
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 74

$ws.Range("A86").Value = 'Kenia'
$ws.Range("B86").Value = 4478
$ws.Range("D86").Value = 1586
$ws.Range("E86").Value = 2771
$ws.Range("H86").Value = 121

# Cuba / Croacia swap places (rows 100-101): Croacia overtakes Cuba
$ws.Range("A100").Value = 'Croacia'
$ws.Range("B100").Value = 2317
$ws.Range("C100").Value = 18
$ws.Range("D100").Value = 2142
$ws.Range("E100").Value = 68
$ws.Range("H100").Value = 107

$ws.Range("A101").Value = 'Cuba'
$ws.Range("B101").Value = 2309
$ws.Range("D101").Value = 2071
$ws.Range("E101").Value = 153
$ws.Range("H101").Value = 85

# Islandia (row 111)
$ws.Range("B111").Value = 1823
$ws.Range("C111").Value = 1
$ws.Range("E111").Value = 8

# Benin jumps ahead of Malaui/Crucero/.../Jamaica (rows 139-147)
$ws.Range("A139").Value = 'Benin'
$ws.Range("B139").Value = 765
$ws.Range("C139").Value = 115
$ws.Range("D139").Value = 253
$ws.Range("E139").Value = 499
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = 13

$ws.Range("A140").Value = 'Malaui'
$ws.Range("B140").Value = 730
$ws.Range("C140").Value = 110
$ws.Range("D140").Value = 258
$ws.Range("E140").Value = 461
$ws.Range("G140").Value = 3
$ws.Range("H140").Value = 11

$ws.Range("A141").Value = 'Crucero'
$ws.Range("B141").Value = 712
$ws.Range("D141").Value = 651
$ws.Range("E141").Value = 48
$ws.Range("H141").Value = 13

$ws.Range("A142").Value = 'Ruanda'
$ws.Range("B142").Value = 702
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 357
$ws.Range("E142").Value = 343
$ws.Range("H142").Value = 2

$ws.Range("A143").Value = 'Santo Tome y Principe'
$ws.Range("B143").Value = 698
$ws.Range("C143").Value = 5
$ws.Range("D143").Value = 203
$ws.Range("E143").Value = 483
$ws.Range("H143").Value = 12

$ws.Range("A144").Value = 'San Marino'
$ws.Range("B144").Value = 696
$ws.Range("D144").Value = 610
$ws.Range("E144").Value = 44
$ws.Range("H144").Value = 42

$ws.Range("A145").Value = 'Mozambique'
$ws.Range("B145").Value = 688
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 177
$ws.Range("E145").Value = 507
$ws.Range("H145").Value = 4

$ws.Range("A146").Value = 'Malta'
$ws.Range("B146").Value = 665
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 616
$ws.Range("E146").Value = 40
$ws.Range("H146").Value = 9

$ws.Range("A147").Value = 'Jamaica'
$ws.Range("B147").Value = 657
$ws.Range("C147").Value = 5
$ws.Range("D147").Value = 462
$ws.Range("E147").Value = 185
$ws.Range("H147").Value = 10

# Dominica / Fiyi swap places (rows 202-203): Fiyi overtakes Dominica
$ws.Range("A202").Value = 'Fiyi'
$ws.Range("A203").Value = 'Dominica'

# Islas Turcas y Caicos / Santa Sede swap places (rows 208-209)
$ws.Range("A208").Value = 'Santa Sede'
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = 'Islas Turcas y Caicos'
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

# Islas Virgenes Britanicas / Papua Nueva Guinea swap places (rows 213-214)
$ws.Range("A213").Value = 'Papua Nueva Guinea'
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = 'Islas Virgenes Britanicas'
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
